# Auto-update draw results: append the 2025-10-24 Pick 3 row to the
# Results sheet.
#
# The source data is plain text (dates/phase codes/time stamps stored as
# strings, not real Excel dates/numbers), matching every other row already
# in the sheet. A leading apostrophe forces Excel to keep the literal text
# instead of auto-converting look-alike values (e.g. "2025-10-24",
# "251024") into a date serial / number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

$ws.Cells.Item($row, 1).Value = "'2025-10-24"
$ws.Cells.Item($row, 2).Value = "'Pick 3"
$ws.Cells.Item($row, 3).Value = "'251024"
$ws.Cells.Item($row, 4).Value = "'5-4-8"
$ws.Cells.Item($row, 5).Value = "'2025-10-24T21:38:03.434+04:00"
